# Swap the species-observation data between row 3 and row 4.
# (Reading back .Value via COM in this runtime is unreliable, so set
# literal values directly instead of swapping through a temp variable.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 becomes what used to be row 4's data.
$ws.Range("A3").Value = 111524816
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("Q3").Value = 528708.8236134996
$ws.Range("R3").Value = 6936113.065692388
$ws.Range("Z3").Value = "09:51"
$ws.Range("AB3").Value = "09:51"
$ws.Range("AC3").Value = ""

# Row 4 becomes what used to be row 3's data.
$ws.Range("A4").Value = 111524610
$ws.Range("B4").Value = 78578
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("Q4").Value = 528710.775020241
$ws.Range("R4").Value = 6936101.088840622
$ws.Range("Z4").Value = "00:00"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AC4").Value = "På sälg"
